# Update the StructureDefinition-Device.xlsx workbook to match the
# published CDA FHIR logical model (patches #241).

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version: 2.0.0-sd-202312-matchbox-patch -> 2.0.0-sd-202406-matchbox-patch
$meta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"

# Date: 2024-03-12T18:28:21+01:00 -> 2024-06-19T17:47:42+02:00
$meta.Range("B8").Value = "2024-06-19T17:47:42+02:00"

# Contact: "No display for ContactDetail" -> full HL7 contact string
$meta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# --- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Device.code binding value set URL
$elements.Range("Z14").Value = "http://hl7.org/cda/stds/core/ValueSet/CDAEntityCode"

# Device.classCode cardinality: Min 1 -> 0 (and matching Base Min).
# These columns store cardinalities as text (shared strings), so a plain
# Value assignment of "0" would be auto-coerced to a number and also
# change the cell's number format/style. Instead, copy an existing cell
# that already holds the text string "0" (e.g. row 3) so the destination
# keeps the same text type and style as its neighbours (F12/AG12 keep
# s="2", t="s").
$elements.Range("F3").Copy() | Out-Null
$elements.Range("F12").PasteSpecial(-4163) | Out-Null   # xlPasteValues

$elements.Range("AG3").Copy() | Out-Null
$elements.Range("AG12").PasteSpecial(-4163) | Out-Null  # xlPasteValues

$excel.CutCopyMode = 0
